$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Remove the "Meta description" paragraph that currently sits right
# after the title heading near the top of the document. It is made up of an
# empty run, a bold run containing "Meta description" and a plain run with
# the rest of the sentence (": Enjoy a thrilling ... Play for free.").
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Meta description*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# Step 2: Insert a brand-new bold paragraph reading "Play Bruce Lee Dragon's
# Tale Slot for Free - Review 2021" right before the last paragraph of the
# document (the one that used to hold the AI image prompt "Please create a
# cartoon-style image ..."). We insert the raw paragraph XML just before the
# paragraph mark of the preceding paragraph ("Not yet optimized for mobile
# devices ...") so that Word creates a clean, independent paragraph - with
# its own leading empty run followed by a bold run - instead of the new text
# inheriting the italic formatting of the paragraph that follows it.
# ---------------------------------------------------------------------------
$prev = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Not yet optimized for mobile devices*") {
        $prev = $p
        break
    }
}

if ($prev -ne $null) {
    $insertPos = $prev.Range.End - 1
    $insertionRange = $d.Range($insertPos, $insertPos)
    $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                  '<w:r/>' +
                  '<w:r><w:rPr><w:b/></w:rPr>' +
                  '<w:t>Play Bruce Lee Dragon''s Tale Slot for Free - Review 2021</w:t>' +
                  '</w:r></w:p>'
    $insertionRange.InsertXML($newParaXml) | Out-Null
}

# ---------------------------------------------------------------------------
# Step 3: Replace the text of the final paragraph (formerly the AI image
# prompt, still carrying italic formatting) with the new meta-description
# sentence, while keeping the existing italic run formatting intact.
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Please create a cartoon-style*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pStart = $target.Range.Start
    $pEnd = $target.Range.End - 1
    $textRange = $d.Range($pStart, $pEnd)
    $textRange.Text = "Enjoy a thrilling and immersive gaming experience with Bruce Lee Dragon's Tale slot game. Explore features, symbols, and betting options in our review. Play for free."
}
